$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04450869273252551
$ws.Range("H2").Value = 5.274608329090507
$ws.Range("I2").Value = 23.52596972994051

$ws.Range("G3").Value = 0.07439599973514201
$ws.Range("H3").Value = 52.76276443197791

$ws.Range("G4").Value = -0.01173324638103145
$ws.Range("H4").Value = -716.2757381652136

$ws.Range("G5").Value = 0.008626083954053733
$ws.Range("H5").Value = 127.2204053174501

$ws.Range("G6").Value = 0.02740633020401429
$ws.Range("H6").Value = -20.94992026615799

$ws.Range("G7").Value = 0.05354401529554811
$ws.Range("H7").Value = 0.6675613464917063

$ws.Range("G8").Value = 0.005046070831298531
$ws.Range("H8").Value = 126.8111744850031

$ws.Range("G9").Value = -0.0252491864677134
$ws.Range("H9").Value = -17.28585039306263

$ws.Range("G10").Value = -0.0604078553441143
$ws.Range("H10").Value = 16.90958921993137

$ws.Range("G11").Value = -0.06430173820599172
$ws.Range("H11").Value = 30.11519363690258

$ws.Range("G12").Value = -0.202818541055177
$ws.Range("H12").Value = 17.03224060698864

$ws.Range("G13").Value = -0.2710174048521587
$ws.Range("H13").Value = 1.381234883905498

$ws.Range("G14").Value = -0.04750611276353835
$ws.Range("H14").Value = -28.05850832063604

$ws.Range("G15").Value = -0.01870957892959942
$ws.Range("H15").Value = 46.19255817689974

$ws.Range("G16").Value = 0.1287942650911535
$ws.Range("H16").Value = 2.782003085849896

$ws.Range("G17").Value = 0.1409216792392777
$ws.Range("H17").Value = 0.4785637837230569

$ws.Range("G18").Value = 0.1232369200468747
$ws.Range("H18").Value = -1.199207321494694

$ws.Range("G19").Value = 0.1457114665572611
$ws.Range("H19").Value = 9.381625369554191

$ws.Range("G20").Value = 0.03507812621554313
$ws.Range("H20").Value = 2.161029429513628

$ws.Range("G21").Value = 0.05928955084630756
$ws.Range("H21").Value = 2.153155049424744

$ws.Range("G22").Value = -0.05919385822591892
$ws.Range("H22").Value = 25.86680943360906

$ws.Range("G23").Value = -0.0575002116133107
$ws.Range("H23").Value = 8.06796238468637

$ws.Range("G24").Value = 0.1117535495340776
$ws.Range("H24").Value = -5.385418355041416

$ws.Range("G25").Value = 0.1454987381037751
$ws.Range("H25").Value = 15.31923915691886

$ws.Range("G26").Value = 0.04753531482270276
$ws.Range("H26").Value = -4.363408697875203

$ws.Range("G27").Value = 0.07536127464732534
$ws.Range("H27").Value = -13.05461858572612

$ws.Range("G28").Value = -0.05096432199675377
$ws.Range("H28").Value = 19.85990383287983

$ws.Range("G29").Value = -0.07231232546395179
$ws.Range("H29").Value = -1.601895123036894

$ws.Range("G30").Value = 0.06111984987813138
$ws.Range("H30").Value = -4.062710750138296

$ws.Range("G31").Value = 0.05278445143313845
$ws.Range("H31").Value = -12.86885188516301

$ws.Range("G32").Value = 0.0783183582711257
$ws.Range("H32").Value = -20.29800250523797

$ws.Range("G33").Value = 0.1113995908337573
$ws.Range("H33").Value = 35.38310753229134

$ws.Range("G34").Value = 0.0107970822474767
$ws.Range("H34").Value = -58.56076471650052

$ws.Range("G35").Value = 0.01474365874763803
$ws.Range("H35").Value = 231.521102023018

$ws.Range("G36").Value = 0.006795065764982753
$ws.Range("H36").Value = 1162.537413544557

$ws.Range("G37").Value = 0.001969108940494451
$ws.Range("H37").Value = 115.684780827067

$ws.Range("G38").Value = 0.108479814514687
$ws.Range("H38").Value = 1.139826803350913

$ws.Range("G39").Value = 0.1062663709023136
$ws.Range("H39").Value = 24.05197270115042

$ws.Range("G40").Value = 0.006067813824378944
$ws.Range("H40").Value = 104.285798728985

$ws.Range("G41").Value = 0.0168698007032865
$ws.Range("H41").Value = 12.48973754961959

$ws.Range("G42").Value = 0.1001843329772523
$ws.Range("H42").Value = -0.74383674139834

$ws.Range("G43").Value = 0.1347134502623401
$ws.Range("H43").Value = 12.12618692825009

$ws.Range("G44").Value = 0.03056699673655603
$ws.Range("H44").Value = -14.34837924854331

$ws.Range("G45").Value = 0.03876509331194158
$ws.Range("H45").Value = 136.8031819633312

$ws.Range("G46").Value = 0.04465517712249295
$ws.Range("H46").Value = 23.22572273432271

$ws.Range("G47").Value = 0.06723162090787108
$ws.Range("H47").Value = 33.28945788473281

$ws.Range("G48").Value = 0.05358553639104247
$ws.Range("H48").Value = 25.2545705817949

$ws.Range("G49").Value = 0.084104939401968
$ws.Range("H49").Value = 21.05592943432601

$ws.Range("G50").Value = 0.0006833431000162872
$ws.Range("H50").Value = -96.04379834448186

$ws.Range("G51").Value = 0.007507848865195604
$ws.Range("H51").Value = -61.43831386461189

$ws.Range("G52").Value = -0.1111569153117103
$ws.Range("H52").Value = -7.376764204933097

$ws.Range("G53").Value = -0.09219073110842281
$ws.Range("H53").Value = 0.1778757263890471

$ws.Range("G54").Value = 0.08826960813722166
$ws.Range("H54").Value = 20.71443637890905

$ws.Range("G55").Value = 0.09580944725647386
$ws.Range("H55").Value = 54.65195182722762

$ws.Range("G56").Value = 0.03161045802110772
$ws.Range("H56").Value = -9.65656680453198

$ws.Range("G57").Value = 0.008109288848317341
$ws.Range("H57").Value = 40.45663722433145

$ws.Range("G58").Value = 0.03468211169995086
$ws.Range("H58").Value = 38.66961233688706

$ws.Range("G59").Value = 0.01641283047586915
$ws.Range("H59").Value = -30.68503794081029

$ws.Range("G60").Value = 0.01602420464276144
$ws.Range("H60").Value = -50.60760852581718

$ws.Range("G61").Value = 0.05526736010205072
$ws.Range("H61").Value = 336.6202395642575

$ws.Range("G62").Value = 0.05619731348457303
$ws.Range("H62").Value = -6.900470141793799

$ws.Range("G63").Value = 0.04321370065392041
$ws.Range("H63").Value = 32.60009966035366

$ws.Range("G64").Value = 0.02820368658616158
$ws.Range("H64").Value = -30.40625619485041

$ws.Range("G65").Value = 0.0620939168405759
$ws.Range("H65").Value = 10.75952369818973

$ws.Range("G66").Value = 0.1045447927626845
$ws.Range("H66").Value = 11.74794865699239

$ws.Range("G67").Value = 0.124949561456755
$ws.Range("H67").Value = 8.231425078169478

$ws.Range("G68").Value = -0.05220498191930791
$ws.Range("H68").Value = -49.79701175277922

$ws.Range("G69").Value = -0.01266969045580872
$ws.Range("H69").Value = 40.29884885302413

$ws.Range("G70").Value = 0.08134777602641
$ws.Range("H70").Value = -12.18581382687166

$ws.Range("G71").Value = 0.09348140181786885
$ws.Range("H71").Value = 2.491651696207769

$ws.Range("G72").Value = -0.03594124441631346
$ws.Range("H72").Value = 35.9119351243362

$ws.Range("G73").Value = -0.06584968684629698
$ws.Range("H73").Value = 10.72752873806446

$ws.Range("G74").Value = 0.09537061795118693
$ws.Range("H74").Value = -4.578624360946026

$ws.Range("G75").Value = 0.120296983610797
$ws.Range("H75").Value = 23.50309595929182

$ws.Range("G76").Value = 0.019185788046214
$ws.Range("H76").Value = -24.96934262598052

$ws.Range("G77").Value = 0.02744894260319163
$ws.Range("H77").Value = 94.54861787077313

$ws.Range("G78").Value = 0.07777422153060559
$ws.Range("H78").Value = 20.99858348028238

$ws.Range("G79").Value = 0.06132087267246576
$ws.Range("H79").Value = -20.06549572049551

$ws.Range("G80").Value = -0.1657044309178647
$ws.Range("H80").Value = -0.0579853342597589

$ws.Range("G81").Value = -0.1398010375478449
$ws.Range("H81").Value = 33.45723402570566

$ws.Range("G82").Value = 0.1156483132080805
$ws.Range("H82").Value = 0.830884051796439

$ws.Range("G83").Value = 0.1832310302351879
$ws.Range("H83").Value = 2.949068704405804

$ws.Range("G84").Value = 0.04242976263781335
$ws.Range("H84").Value = 77.99474503204689

$ws.Range("G85").Value = 0.03225304838701491
$ws.Range("H85").Value = -47.62069685184704
